$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "(according to the population census data)" subtitle row.
$ws.Rows.Item(2).Delete()

# Remove the 1989 and 2002 columns, keeping only the 2014 figures.
$ws.Range("B1:C1").EntireColumn.Delete()

# Match the refreshed row heights used across the re-exported sheet.
$ws.Range("A1:B5").EntireRow.RowHeight = 20.1
